$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename the smoke test identifier in A6
$ws.Range("A6").Value = "102_AutomobileInsurance_001_SmokeTest_FillPage"

# Update the selection on the sheet (active cell moved from J1 to I14)
$ws.Range("I14").Select()
